# Add the new follow-up row/question to the sheet: a new cell A9 holding the
# text "The last time I did not complete it" (this introduces a new shared
# string and extends the used range from A1:D7 to A1:D9), then move the
# active selection onto the newly added cell (A9), matching where the user
# left off editing.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value = "The last time I did not complete it"
$null = $ws.Range("A9").Select()
